$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is added for Femacal de La Calera - Mango.
# It belongs at row 419 (in date order among the existing rows), so the
# existing rows 419:447 shift down to 420:448.
$ws.Rows("419:419").Insert()

# Populate the newly inserted row 419 with the new record. It shares the
# same market/product/category/quality/origin/unit metadata as the row
# that used to sit there (now at row 420), only the date, volume, prices
# and $/Kg are new.
$ws.Range("A419").Value = 3
$ws.Range("B419").Value = "Femacal de La Calera"
$ws.Range("C419").Value = "Coquimbo"
$ws.Range("D419").Value = 44783
$ws.Range("E419").Value = 5
$ws.Range("F419").Value = "Fruta"
$ws.Range("G419").Value = 100108
$ws.Range("H419").Value = "Tropicales y subtropicales"
$ws.Range("I419").Value = 100108002
$ws.Range("J419").Value = "Mango"
$ws.Range("K419").Value = "Sin especificar"
$ws.Range("L419").Value = "Primera"
$ws.Range("M419").Value = 228
$ws.Range("N419").Value = 10000
$ws.Range("O419").Value = 10000
$ws.Range("P419").Value = 10000
$ws.Range("Q419").Value = "$/bandeja 4 kilos"
$ws.Range("R419").Value = "Brasil"
$ws.Range("S419").Value = 2500
$ws.Range("T419").Value = 4
